# Update gh-pages output data for 江西-漫展信息.xlsx
# Applies refreshed "想去人数" (want-to-go count) figures and marks one
# event as sold out ("已售罄") on both the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# Map of row -> new F value (想去人数) for the affected rows.
$fUpdates = @{
    3  = 8564
    4  = 1525
    6  = 395
    7  = 265
    9  = 33
    13 = 1260
    14 = 264
    15 = 83
    16 = 146
    17 = 101
    18 = 134
    19 = 82
    20 = 122
    21 = 108
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    # Row 3's lowest ticket price column is now sold out, switch it from a
    # numeric price to the text marker used elsewhere in the sheet.
    $ws.Range("G3").Value = "已售罄"
}
